$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new day's fuel-price reading (2025-04-24, serial 45771) needs to go on top
# of the table; every existing data row (2..30) shifts down by one (3..31).
$ws.Rows.Item(2).Insert()

# Rows.Insert() pulls formatting from the row above (the bold header row), so
# fix the new row's formatting by copying it from the row right below instead
# (which still carries the normal data-row formatting: date format in col A,
# default/general format in cols B and C).
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("B3").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C2").PasteSpecial(-4122)

# Now write the new day's values into row 2.
$ws.Cells.Item(2, 1).Value = 45771
$ws.Cells.Item(2, 2).Value = 700.96
$ws.Cells.Item(2, 3).Value = 753.256
